$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the cells whose new value is a plain decimal number as Text,
# so Excel stores the literal string (e.g. "1.00") instead of silently
# re-interpreting it as a number and dropping trailing/insignificant zeros.
$textCells = @("D4","D5","D6","D10","D12","D16","D19","D21","D22","D24","D26","D27","D30","D31","D32","D33","D34","D35","D36","D38","D39","D41","D44","D45","D46","D47","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value = "64.358.47"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.456.34"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "574.60"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "160.12"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "3.458.20"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("E9").Value = "  +10.79%  "
$ws.Range("D10").Value = "7.36"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "4.044.02"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("E15").Value = "  +5.59%  "
$ws.Range("D16").Value = "28.84"
$ws.Range("E16").Value = "  +6.39%  "
$ws.Range("D17").Value = "64.421.52"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "3.444.34"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +3.61%  "
$ws.Range("D21").Value = "386.64"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "8.23"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("E23").Value = "  +1.85%  "
$ws.Range("D24").Value = "73.14"
$ws.Range("E24").Value = "  +3.24%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  +17.89%  "
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").Value = "  +10.55%  "
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  +8.89%  "
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "6.61"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "23.76"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "7.08"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "160.49"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "0.0778"
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "27.48"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "2.914.14"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("D44").Value = "4.49"
$ws.Range("E44").Value = "  +4.33%  "
$ws.Range("D45").Value = "42.42"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("D46").Value = "0.770"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "23.83"
$ws.Range("E47").Value = "  +7.87%  "
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").Value = "2.20"
$ws.Range("E49").Value = "  +15.05%  "
$ws.Range("E51").Value = "  +3.93%  "
